$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1933449667456422000
$ws.Range("B2").Value = "@LeticiaFrost0 @_pao_com_banana Fui lembramdo por lele da silva yayyy https://t.co/ReFEDaDEEp"
$ws.Range("C2").Value = "June 13, 2025 at 09:01 AM"
$ws.Range("D2").Value = " pao com banana fui lembramdo por lele da silva yayyy "
$ws.Range("E2").Value = "es"

$ws.Range("A3").Value = 1933261971622613000
$ws.Range("B3").Value = "@p2rcys fun fact no inverno tem aumento de casos de infeccao urinaria"
$ws.Range("C3").Value = "June 12, 2025 at 08:35 PM"
$ws.Range("D3").Value = " fun fact no inverno tem aumento de casos de infeccao urinaria"
$ws.Range("E3").Value = "es"

$ws.Range("A4").Value = 1932477783604712000
$ws.Range("B4").Value = "@ybernxd infeccao urinaria"
$ws.Range("C4").Value = "June 10, 2025 at 04:39 PM"
$ws.Range("D4").Value = " infeccao urinaria"
$ws.Range("E4").Value = "es"

$ws.Range("A5").Value = 1931740333252944000
$ws.Range("B5").Value = "@psyllobor4 lembrando q sodio da infeccao urinaria"
$ws.Range("C5").Value = "June 8, 2025 at 03:49 PM"
$ws.Range("D5").Value = " lembrando q sodio da infeccao urinaria"
$ws.Range("E5").Value = "es"

$ws.Range("A6").Value = 1931476868613477000
$ws.Range("B6").Value = "@tudojaywon a infeccao urinaria babado"
$ws.Range("C6").Value = "June 7, 2025 at 10:22 PM"
$ws.Range("D6").Value = " a infeccao urinaria babado"
$ws.Range("E6").Value = "es"

$ws.Range("A7").Value = 1903095762231182000
$ws.Range("B7").Value = "infeccao urinaria eh para poucos https://t.co/2rRhtIWgNv"
$ws.Range("C7").Value = "March 21, 2025 at 02:46 PM"
$ws.Range("D7").Value = "feccao urinaria eh para poucos "
$ws.Range("E7").Value = "es"

$ws.Range("A8").Value = 1884593735926797000
$ws.Range("B8").Value = "Acho q to c infeccao urinaria"
$ws.Range("C8").Value = "January 29, 2025 at 01:25 PM"
$ws.Range("D8").Value = "acho q to c infeccao urinaria"
$ws.Range("E8").Value = "es"

$ws.Range("A9").Value = 1884359698993800000
$ws.Range("B9").Value = "acho que to infeccao urinaria"
$ws.Range("C9").Value = "January 28, 2025 at 09:55 PM"
$ws.Range("D9").Value = "acho que to infeccao urinaria"
$ws.Range("E9").Value = "es"

$ws.Range("A10").Value = 1883974742761431000
$ws.Range("B10").Value = "@conexaocec @venecasagrande DUROS"
$ws.Range("C10").Value = "January 27, 2025 at 08:25 PM"
$ws.Range("D10").Value = " duros"
$ws.Range("E10").Value = "es"

$ws.Range("A11:E13").EntireRow.Delete()

